# Freelancer k=0.1 results workbook: add summary rows (Average/Worst of
# SW(S*)/SW(OPT) and SC(S*)/SC(OPT)) below the per-instance data, plus the
# stray AVERAGE() helper formula that sits in the now-blank row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J12: leftover helper formula in the (otherwise blank) row under the data.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# --- New summary rows under the data table.
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold 12pt, vertically centered -- set once on B14, then fan the exact same
# style out to B15:B17 via a format-only paste so they all share one cellXf
# instead of each deriving (and leaving behind) their own intermediate style.
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

# Row heights match the taller summary rows.
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Page setup (A4 portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection mirrors the state the workbook was saved in.
$ws.Range("A14:B17").Select()
